# Update gh-pages to output generated at 456a3b4
# This script applies numeric "tickets sold" bumps across the four
# worksheets of the 上海-漫展信息 workbook, switches one cell from a
# sold-count number to a "已售罄" (sold out) text flag, and refreshes the
# "全部类型" (all types) sheet so that a newly-added 本地生活/展览 entry
# (洛天依歌行宇宙·无限遨游 沉浸式体验展) is reflected there as well,
# cascading the rows that were already present.

$wb = $excel.ActiveWorkbook

# Worksheets, identified by name to be robust against ordering.
$wsExpo    = $wb.Worksheets.Item("展览")     # sheet1
$wsShow    = $wb.Worksheets.Item("演出")     # sheet2
$wsLocal   = $wb.Worksheets.Item("本地生活") # sheet3
$wsAll     = $wb.Worksheets.Item("全部类型") # sheet4

# ---------------------------------------------------------------
# Sheet "展览" (sheet1): simple numeric refreshes in column F.
# ---------------------------------------------------------------
$wsExpo.Range("F6").Value  = 55
$wsExpo.Range("F8").Value  = 329
$wsExpo.Range("F9").Value  = 2180
$wsExpo.Range("F12").Value = 842
$wsExpo.Range("F15").Value = 1460
$wsExpo.Range("F16").Value = 694
$wsExpo.Range("F17").Value = 1680
$wsExpo.Range("F18").Value = 38
$wsExpo.Range("F19").Value = 343
$wsExpo.Range("F23").Value = 2611

# ---------------------------------------------------------------
# Sheet "演出" (sheet2): simple numeric refreshes in column F.
# ---------------------------------------------------------------
$wsShow.Range("F19").Value = 150
$wsShow.Range("F33").Value = 60
$wsShow.Range("F38").Value = 334
$wsShow.Range("F43").Value = 74

# ---------------------------------------------------------------
# Sheet "本地生活" (sheet3): numeric refreshes plus one cell that
# switches from a numeric remaining-stock count to a "sold out" flag.
# ---------------------------------------------------------------
$wsLocal.Range("F4").Value  = 2479
$wsLocal.Range("G6").Value  = "已售罄"
$wsLocal.Range("F13").Value = 2759
$wsLocal.Range("F14").Value = 351
$wsLocal.Range("F15").Value = 660

# ---------------------------------------------------------------
# Sheet "全部类型" (sheet4): numeric refresh of the row mirroring
# "本地生活"!F4, then a cascade of rows 4-6 down by one (to make room
# for the already-existing-elsewhere "洛天依歌行宇宙" entry that is
# newly surfaced on this aggregate sheet), followed by more numeric
# refreshes on the unaffected rows further down.
# ---------------------------------------------------------------
$wsAll.Range("F2").Value = 2479

# Column B holds plain "YYYY-MM-DD" text, not a real date, throughout this
# workbook. Force the B4:B7 range through a Text number format first so
# Excel does not auto-convert the literals into date serials, then restore
# the default "Normal" style afterwards so no stray formatting is left on
# the cells.
$wsAll.Range("B4:B7").NumberFormat = "@"

# Row 4 <- old Row 5 (上海·日漫咖啡体验)
$wsAll.Range("B4").Value = "2024-09-09"
$wsAll.Range("C4").Value = "上海·日漫咖啡体验"
$wsAll.Range("D4").Value = "虹桥路1438号高岛屋百货6楼 Oasis漫画喫茶"
$wsAll.Range("E4").Value = "2024.09.09 10:00-12.31 22:00"
$wsAll.Range("F4").Value = 131
$wsAll.Range("G4").Value = 60
$wsAll.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=91993"
$wsAll.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202409/IV5rInWT1725347808557.jpeg"

# Row 5 <- old Row 6 (上海·迷你四驱车赛场)
$wsAll.Range("B5").Value = "2024-09-10"
$wsAll.Range("C5").Value = "上海·迷你四驱车赛场"
$wsAll.Range("D5").Value = "虹桥路1438号高岛屋百货6楼 Oasis漫画喫茶"
$wsAll.Range("E5").Value = "2024.09.10 10:00-12.31 22:00"
$wsAll.Range("F5").Value = 6
$wsAll.Range("G5").Value = 48
$wsAll.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=92042"
$wsAll.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202409/LzFT5TMO1725348229429.png"

# Row 6 <- old Row 7 (上海·吉卜力工作室物语-沉浸式艺术展全球首站（9月-10月）)
$wsAll.Range("B6").Value = "2024-09-14"
$wsAll.Range("C6").Value = "上海·吉卜力工作室物语-沉浸式艺术展全球首站（9月-10月）"
$wsAll.Range("D6").Value = "龙台路10号2F 上海国际传媒港艺术中心"
$wsAll.Range("E6").Value = "2024.09.14 10:00-10.31 20:00"
$wsAll.Range("F6").Value = 161
$wsAll.Range("G6").Value = 88
$wsAll.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=91856"
$wsAll.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202409/wSR0yFfg1725432304586.jpeg"

# Row 7 <- brand new entry (上海 洛天依歌行宇宙·无限遨游 沉浸式体验展)
$wsAll.Range("B7").Value = "2024-09-15"
$wsAll.Range("C7").Value = "上海 洛天依歌行宇宙·无限遨游 沉浸式体验展"
$wsAll.Range("D7").Value = "中山北路3300号 上海月星环球港"
$wsAll.Range("E7").Value = "2024.09.15 10:00-10.31 22:00"
$wsAll.Range("F7").Value = 2407
$wsAll.Range("G7").Value = 138
$wsAll.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=91175"
$wsAll.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202408/ei9COXS41724405861343.jpeg"

$wsAll.Range("B4:B7").Style = "Normal"

# Remaining untouched rows: numeric refreshes only.
$wsAll.Range("F8").Value  = 2759
$wsAll.Range("F9").Value  = 351
$wsAll.Range("F11").Value = 660
$wsAll.Range("F16").Value = 55
$wsAll.Range("F18").Value = 329
$wsAll.Range("F21").Value = 842
$wsAll.Range("F28").Value = 694
$wsAll.Range("F31").Value = 1680
$wsAll.Range("F32").Value = 343
$wsAll.Range("F42").Value = 334
$wsAll.Range("F43").Value = 2611
$wsAll.Range("F46").Value = 74
